# Fixes indentation (fill in médiation sheet data) + switches active tab
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range('A1').Value = 'page'
$ws2.Range('A2').Value = 'code'
$ws2.Range('B2').Value = 'THF'
$ws2.Range('C2').Value = 'EG'
$ws2.Range('D2').Value = 'AMAT'
$ws2.Range('E2').Value = 'SCO'
$ws2.Range('F2').Value = 'TEAM'
$ws2.Range('G2').Value = 'MES'
$ws2.Range('A3').Value = 'titre'
$ws2.Range('B3').Value = 'Théâtre-Forum'
$ws2.Range('C3').Value = 'Escape-Game'
$ws2.Range('D3').Value = 'accompagnement amateurs'
$ws2.Range('E3').Value = 'interventions scolaires'
$ws2.Range('F3').Value = 'team Building'
$ws2.Range('G3').Value = 'Intervention sur mesure'
$ws2.Range('A4').Value = 'type'
$ws2.Range('B4').Value = 'Animation-débat'
$ws2.Range('C4').Value = 'Jeu collectif'
$ws2.Range('D4').Value = 'Mise-en-scène et direction d''acteurs'
$ws2.Range('E4').Value = 'du CP à la terminale'
$ws2.Range('F4').Value = 'Activités coopératives'
$ws2.Range('G4').Value = 'Aucune limite à vos envies'
$ws2.Range('A5').Value = 'durée'
$ws2.Range('B5').Value = 'de 45 à 100 minutes'
$ws2.Range('C5').Value = 'de 55 à 90 minutes'
$ws2.Range('D5').Value = 'accompagnement ponctuel ou hebdomadaire'
$ws2.Range('E5').Value = 'des séquences construites avec vous'
$ws2.Range('F5').Value = 'de 1 à 4h'
$ws2.Range('G5').Value = 'à votre convenance'
$ws2.Range('A6').Value = 'date création'
$ws2.Range('B6').Value = 'depuis 2012'
$ws2.Range('C6').Value = 'depuis 2019'
$ws2.Range('D6').Value = 'depuis 2010'
$ws2.Range('E6').Value = 'depuis 2010'
$ws2.Range('F6').Value = 'depuis 2015'
$ws2.Range('G6').Value = 'depuis toujours !'
$ws2.Range('A7').Value = 'accroche'
$ws2.Range('B7').Value = 'Une autre façon de créer le débat'
$ws2.Range('C7').Value = 'Un savoir-faire unique !'
$ws2.Range('D7').Value = 'Avancer ensemble'
$ws2.Range('E7').Value = 'Au service du savoir'
$ws2.Range('F7').Value = 'La cohésion avant tout'
$ws2.Range('G7').Value = 'Notre credo : vous proposer mieux que ce que vous attendiez !'
$ws2.Range('A8').Value = 'qu''est-ce que c''est'
$ws2.Range('B8').Value = 'Le théâtre-forum est un outil puissant pour aborder des problématiques d''une manière nouvelle. Quelques semaines avant la prestation, nous rencontrons les organisateurs pour déterminer vos problématiques. Ensuite, une saynète est écrite. Elle est jouée devant les spectateurs. Cette scène a a particularité d''être très réaliste et de mal se terminer. petit à petit, les spectateurs font des propositions pour aider à améliorer la situation. Ils viennent les tester sur scène avec les comédiens, et on avance ainsi pas à pas.'
$ws2.Range('C8').Value = 'Vous connaissez probablement l''escape game traditionnel : une équipe de 2 à 6 personnes dispose d''une heure pour réaliser une mission… Nous sommes allé plus loin : tout d''abord, nous avons créé des enquêtes d''extérieur ; les fondamentaux de l''escape game sont toujours là, mais votre terrain de jeu est une ville ! Nous avons aussi créé des parties d''escape Game "grand format" qui font jouer jusqu''à 35 joueurs simultanément ! Des formules à découvrir !'
$ws2.Range('D8').Value = 'L''accompagnement des comédiens amateurs fait partie de notre ADN. Nous proposons plusieurs types d''accompagnement principalement aurpès de la compagnie VOLT qui se fait le relais entre les groupements de comédiens et nous. Nous pouvons proposer des stages plus ou moins longs autour de techniques théâtrales (mime, improvisation, pour adultes ou pour enfants, etc.). Nous pouvons aussi accompagner un projet de pièce de théâtre, en nous positionnant en fonction des attentes des comédiens : la mise-en-scène, la direction d''acteurs, un simple regard extérieur, des apports techniques, etc. '
$ws2.Range('E8').Value = 'Notre panel d''interventions est très large. Cela va de l''intervention ponctuelle, à l''animation d''ateliers théâtre ou improvisation hebdomadaires. Chaque année, nous intervenons auprès de 5 ou 6 classes dans des projets aussi divers que "les lectures thétralisées", des "ateliers d''impro", "l''accompagnement à l''écritreu dramatique", "la mise en scène", etc.  Tout cela sans compter, évidemment les interventions en théâtre-forum ou les spectacles de fin d''année.'
$ws2.Range('F8').Value = 'Les entreprises cherchent parfois à proposer un temps différents à leurs salariés. Les objectifs principaux : créer une dynamique de groupe et améliorer la cohésion de l''équipe. Nous proposons plusieurs choses en team-building. Il s''agit principalement de jeux coopératifs qui s''appuient sur notre expérience en improvisation théâtrale (la culture du "oui") couplés parfois avec une expérience d''escape-game grand format.'
$ws2.Range('G8').Value = 'Ici : seule votre imagination ou vos envies sont les limites ! Nous sommes à votre écoute et vous ferons des propositions en lien avec vos objectifs. Ne vous censurez pas : la folie est notre moteur !'
$ws2.Range('A9').Value = 'pour qui ?'
$ws2.Range('B9').Value = 'Nous intervenons pour les associsations d''aide aux personnes ou les établissements scolaires.  Pour une meilleure personnalisation, nous pouvons créer la ou les saynètes en fonction de vos problématiques.Pour les petits budgets, nous disposons d''un grand stock de saynètes déjà écrites.'
$ws2.Range('C9').Value = 'Nous conseillons aux particuliers de passer directement par notre partenaire : l''agence E.SCAPE de Valréas. Pour les entreprises, les grands groupes, ou les événements ponctuels : contactez-nous car nous avons probablement des choses à vous proposer !'
$ws2.Range('D9').Value = 'Nous n''intervenons pas directement auprès des particuliers. Nous vous invitons à vous rapprocher de notre partenaire la compagnie VOLT de Nyons, ou à nous contacter via une structure de théâtre amateur.'
$ws2.Range('E9').Value = 'Nous intervenons pour toutes les établissements scolaires sans exception, ainsi que pour les établissements privés, sous-contrat ou hors-contrat.'
$ws2.Range('F9').Value = 'Pour les entreprises de plus de 10 salariés ou les associations qui souhaitent créer une émulations entre membres.'
$ws2.Range('G9').Value = 'Toute structure qui cherche à créer un événement innovant.'
$ws2.Range('A10').Value = 'notre expérience'
$ws2.Range('B10').Value = 'Nous disposons de dizaines de saynètes déjà écrites, mais les écrivons pour les demandes spécifiques. Nos prestations tournent depuis 10 ans, avec un succès et des effets qui ne se démentent pas.'
$ws2.Range('C10').Value = 'Nous avons créer l''Agence E.SCAPE, devenue autonome en 2019. Elle propose nos enquêtes, qui sont au nombre de 3 (deux en intérieur, une en extérieur), et qui font le bonheur de tous les joueurs ! Notre particularité est d''apporter un soin tout particulier à l''encadrement et à l''immersion des joueurs. Nous avons à coeur de faire en sorte que l''escape Game soit une sorte de pièce de théâtre dont vous êtes le héros !'
$ws2.Range('D10').Value = 'Nous avons accompagné plusieurs centaines d''amateurs, adultes ou enfants et monté plusieurs dizaines de pièces de théâtre pour amateurs.'
$ws2.Range('E10').Value = 'Nous comptons dans nos intervenants d''anciens enseignants titulaires, devenus comédiens, qui ont une parfaite connaissance du milieu scolaire. Forts de cet atout, nous savons nous adpater à vos contraintes et veillons à ce que nos interventions constituent un plus dans les apprentissages. Parmi les expériences marquantes et originales, mentionnons l''accompagnement d''un groupe d''élèves à l''Ecriture de monologues théâtraux, qu''ils ont mis en scène eux-mêmes afin que notre comédien les joue !!! Nous sommes aussi intervenus afin de proposer une épreuve d''écriture-interprétation de théâtre-forum pour le BAC 2016.'
$ws2.Range('F10').Value = 'L''improvisation théâtrale est un outil puissant pour apprendre à un groupe à coopérer. Et si cela fonctionne dans les établissements scolaires avec des élèves en ruptures, il n''y a aucun raison que ça ne fonctionne pas en entreprise ! Notre positionnement ludique contribuera à rendre l''événement festif, et augmentera l''impact de l''intervention sur votre équipe !'
$ws2.Range('G10').Value = 'De nombreuses structures travaillent avec nous chaque année dans des projets originaux. Nous avons ainsi participer à l''écriture d''une BD pour les points info-énergie, nous avons organisé "un débat sans thème" pour des fédérations de centres sociaux, créé un escape-game pour 50 personnes en un mois... On n''attend que vous pour nous lancer dans un projet fou !'

# Row heights (auto-sized after content/formatting changes)
$ws2.Rows.Item(7).RowHeight = 30
$ws2.Rows.Item(8).RowHeight = 168.75
$ws2.Rows.Item(9).RowHeight = 90
$ws2.Rows.Item(10).RowHeight = 195

# E8 needs wrap + vertical-center alignment (style index 2)
$ws2.Range('E8').WrapText = $true
$ws2.Range('E8').VerticalAlignment = -4108

# Remove leftover "tout-public" placeholder cells from header row (B1:E1)
$ws2.Range('B1:E1').Clear()

# Ludo Quizz "pour qui" cell on tout-public sheet gets real copy instead of a blank placeholder
$ws1.Range('H9').Value = 'Le Ludo quiz est adaptable au possible ! Il nécessite juste une arrivée électrique et, lorsqu''il se jout en extérieur, un endroit ombragé pour la visibilité de l''écran. En version "quiz de rue", il donnera du peps à vos animations thématiques. En version "tournoi de salle", il rythmera votre événement (repas, fête, etc.) et contribuera à créer une ambiance festive !'

# Print setup restored on médiation sheet
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# médiation sheet becomes the active tab/view; selection moves to H7, scrolled to show column D onward
$ws2.Activate()
$ws2.Range('D1').Select()
$ws2.Range('H7').Select()
